$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Cell='D2'; Value='276.40'; Numeric=$true}
    @{Cell='G2'; Value='3'; Numeric=$true}
    @{Cell='D3'; Value='23.15'; Numeric=$true}
    @{Cell='G3'; Value='3'; Numeric=$true}
    @{Cell='D4'; Value='6.461'; Numeric=$true}
    @{Cell='G4'; Value='3'; Numeric=$true}
    @{Cell='D5'; Value='0.06295'; Numeric=$true}
    @{Cell='G5'; Value='3'; Numeric=$true}
    @{Cell='D6'; Value='3.661'; Numeric=$true}
    @{Cell='G6'; Value='3'; Numeric=$true}
    @{Cell='D7'; Value='6.681'; Numeric=$true}
    @{Cell='G7'; Value='3'; Numeric=$true}
    @{Cell='D8'; Value='1.381'; Numeric=$true}
    @{Cell='G8'; Value='3'; Numeric=$true}
    @{Cell='D9'; Value='0.8356'; Numeric=$true}
    @{Cell='G9'; Value='3'; Numeric=$true}
    @{Cell='D10'; Value='0.01383'; Numeric=$true}
    @{Cell='G10'; Value='3'; Numeric=$true}
    @{Cell='D11'; Value='0.1606'; Numeric=$true}
    @{Cell='G11'; Value='3'; Numeric=$true}
    @{Cell='D12'; Value='0.08292'; Numeric=$true}
    @{Cell='G12'; Value='3'; Numeric=$true}
    @{Cell='D13'; Value='0.03441'; Numeric=$true}
    @{Cell='G13'; Value='3'; Numeric=$true}
    @{Cell='D14'; Value='0.03109'; Numeric=$true}
    @{Cell='G14'; Value='3'; Numeric=$true}
    @{Cell='D15'; Value='0.09307'; Numeric=$true}
    @{Cell='G15'; Value='3'; Numeric=$true}
    @{Cell='D16'; Value='3.841'; Numeric=$true}
    @{Cell='G16'; Value='3'; Numeric=$true}
    @{Cell='D17'; Value='0.001647'; Numeric=$true}
    @{Cell='G17'; Value='3'; Numeric=$true}
    @{Cell='D18'; Value='0.04769'; Numeric=$true}
    @{Cell='G18'; Value='3'; Numeric=$true}
    @{Cell='D19'; Value='0.006405'; Numeric=$true}
    @{Cell='G19'; Value='3'; Numeric=$true}
    @{Cell='D20'; Value='0.005689'; Numeric=$true}
    @{Cell='G20'; Value='3'; Numeric=$true}
    @{Cell='D21'; Value='0.001089'; Numeric=$true}
    @{Cell='G21'; Value='3'; Numeric=$true}
    @{Cell='D22'; Value='0.0001499'; Numeric=$true}
    @{Cell='G22'; Value='3'; Numeric=$true}
    @{Cell='D23'; Value='3.714'; Numeric=$true}
    @{Cell='G23'; Value='3'; Numeric=$true}
    @{Cell='G24'; Value='3'; Numeric=$true}
    @{Cell='D25'; Value='0.3346'; Numeric=$true}
    @{Cell='G25'; Value='3'; Numeric=$true}
    @{Cell='D26'; Value='0.1261'; Numeric=$true}
    @{Cell='G26'; Value='3'; Numeric=$true}
    @{Cell='D27'; Value='0.0002678'; Numeric=$true}
    @{Cell='G27'; Value='3'; Numeric=$true}
    @{Cell='G28'; Value='3'; Numeric=$true}
    @{Cell='G29'; Value='3'; Numeric=$true}
    @{Cell='G30'; Value='3'; Numeric=$true}
    @{Cell='G31'; Value='3'; Numeric=$true}
    @{Cell='G32'; Value='3'; Numeric=$true}
    @{Cell='G33'; Value='3'; Numeric=$true}
    @{Cell='G34'; Value='3'; Numeric=$true}
    @{Cell='G35'; Value='3'; Numeric=$true}
    @{Cell='G36'; Value='3'; Numeric=$true}
    @{Cell='G37'; Value='3'; Numeric=$true}
    @{Cell='G38'; Value='3'; Numeric=$true}
    @{Cell='G39'; Value='3'; Numeric=$true}
    @{Cell='D40'; Value='0.04733'; Numeric=$true}
    @{Cell='G40'; Value='3'; Numeric=$true}
    @{Cell='D41'; Value='0.007067'; Numeric=$true}
    @{Cell='G41'; Value='3'; Numeric=$true}
    @{Cell='D42'; Value='0.1164'; Numeric=$true}
    @{Cell='G42'; Value='3'; Numeric=$true}
    @{Cell='D43'; Value='0.003501'; Numeric=$true}
    @{Cell='E43'; Value='42CEJICEJIWorstin24h'; Numeric=$false}
    @{Cell='G43'; Value='3'; Numeric=$true}
    @{Cell='D44'; Value='0.01217'; Numeric=$true}
    @{Cell='G44'; Value='3'; Numeric=$true}
    @{Cell='D45'; Value='0.00006248'; Numeric=$true}
    @{Cell='G45'; Value='3'; Numeric=$true}
    @{Cell='G46'; Value='3'; Numeric=$true}
    @{Cell='D47'; Value='0.00000000749'; Numeric=$true}
    @{Cell='G47'; Value='3'; Numeric=$true}
    @{Cell='D48'; Value='0.7959'; Numeric=$true}
    @{Cell='E48'; Value='47CoinbaseStockTokenCOIN'; Numeric=$false}
    @{Cell='G48'; Value='3'; Numeric=$true}
    @{Cell='B49'; Value='BOLO'; Numeric=$false}
    @{Cell='C49'; Value='https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'; Numeric=$false}
    @{Cell='D49'; Value='0.002122'; Numeric=$true}
    @{Cell='E49'; Value='48BOLOBOLOBestin24h'; Numeric=$false}
    @{Cell='G49'; Value='3'; Numeric=$true}
    @{Cell='B50'; Value='CryptobidCoin'; Numeric=$false}
    @{Cell='C50'; Value='https://coinranking.com/coin/h39bvStAP+cryptobidcoin-cbc'; Numeric=$false}
    @{Cell='D50'; Value='0.00002298'; Numeric=$true}
    @{Cell='E50'; Value='49CryptobidCoinCBC'; Numeric=$false}
    @{Cell='G50'; Value='3'; Numeric=$true}
    @{Cell='D51'; Value='0.01239'; Numeric=$true}
    @{Cell='G51'; Value='3'; Numeric=$true}
)

foreach ($item in $changes) {
    if ($item.Numeric) {
        $ws.Range($item.Cell).Value = "'" + $item.Value
        $ws.Range($item.Cell).Style = "Normal"
    } else {
        $ws.Range($item.Cell).Value = $item.Value
    }
}
